$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.715.29"
$ws.Range("E2").Value = "  +1.34%  "

$ws.Range("D3").Value = "3.467.81"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.28%  "

$ws.Range("D5").Value = "'414.05"
$ws.Range("E5").Value = "  +1.42%  "

$ws.Range("D6").Value = "'129.83"
$ws.Range("E6").Value = "  +0.73%  "

$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  -0.28%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.727"
$ws.Range("E9").Value = "  -0.52%  "

$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  +9.83%  "

$ws.Range("D11").Value = "'42.54"
$ws.Range("E11").Value = "  -0.18%  "

$ws.Range("D12").Value = "'9.71"
$ws.Range("E12").Value = "  +6.38%  "

$ws.Range("D13").Value = "'0.0000224"
$ws.Range("E13").Value = "  +3.68%  "

$ws.Range("D14").Value = "4.011.31"
$ws.Range("E14").Value = "  +1.09%  "

$ws.Range("D15").Value = "'0.141"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("D16").Value = "'20.60"
$ws.Range("E16").Value = "  -3.03%  "

$ws.Range("D17").Value = "3.441.70"
$ws.Range("E17").Value = "  +1.13%  "

$ws.Range("D18").Value = "'12.67"
$ws.Range("E18").Value = "  +1.48%  "

$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").Value = "62.628.80"
$ws.Range("E20").Value = "  +1.20%  "

$ws.Range("D21").Value = "'464.72"
$ws.Range("E21").Value = "  +2.39%  "

$ws.Range("D22").Value = "'90.73"
$ws.Range("E22").Value = "  -0.74%  "

$ws.Range("D23").Value = "'3.28"
$ws.Range("E23").Value = "  +2.32%  "

$ws.Range("D24").Value = "'13.30"
$ws.Range("E24").Value = "  +1.94%  "

$ws.Range("D25").Value = "'10.60"
$ws.Range("E25").Value = "  +15.74%  "

$ws.Range("E26").Value = "  +1.24%  "

$ws.Range("D27").Value = "'33.40"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").Value = "'7.55"
$ws.Range("E29").Value = "  -0.46%  "

$ws.Range("D30").Value = "'12.06"
$ws.Range("E30").Value = "  +0.37%  "

$ws.Range("D31").Value = "'2.67"
$ws.Range("E31").Value = "  -2.91%  "

$ws.Range("E32").Value = "  -1.68%  "

$ws.Range("E33").Value = "  -1.04%  "

$ws.Range("D34").Value = "'40.80"
$ws.Range("E34").Value = "  -4.47%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "'58.50"
$ws.Range("E36").Value = "  +8.47%  "

$ws.Range("D37").Value = "'0.0492"
$ws.Range("E37").Value = "  -1.72%  "

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'3.09"
$ws.Range("E38").Value = "  +5.23%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.06%  "

$ws.Range("D40").Value = "'148.85"
$ws.Range("E40").Value = "  +4.57%  "

$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").Value = "'0.322"
$ws.Range("E41").Value = "  +1.76%  "

$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'2.71"
$ws.Range("E42").Value = "  +6.23%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("E44").Value = "  -1.24%  "

$ws.Range("E45").Value = "  +3.78%  "

$ws.Range("D46").Value = "'2.07"
$ws.Range("E46").Value = "  +3.81%  "

$ws.Range("E47").Value = "  +12.88%  "

$ws.Range("D48").Value = "'" + "0.0" + [char]0x2083 + "0570"
$ws.Range("E48").Value = "  +35.43%  "

$ws.Range("D49").Value = "'16.41"
$ws.Range("E49").Value = "  -0.87%  "

$ws.Range("D50").Value = "'22.27"
$ws.Range("E50").Value = "  -0.07%  "

$ws.Range("E51").Value = "  -1.95%  "
